$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18: You Grow, Girl | Growth Formula Beta
$ws.Range("H18").Value = 15765.934
$ws.Range("I18").Value = 13187.5
$ws.Range("J18").Value = 18712.715
$ws.Range("K18").Value = 13187.5
$ws.Range("L18").Value = 18712.715
$ws.Range("M18").Value = -12903.5
$ws.Range("N18").Value = -19280.715

# Row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws.Range("H132").Value = 1328.5
$ws.Range("I132").Value = 1276.8462
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 3830.5386
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -1300.5386
$ws.Range("N132").Value = -11060


$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Range("H61").Value = 3522
$ws.Range("I61").Value = 2044
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 2044
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -1832
$ws.Range("N61").Value = -5424

# Row 74: As the Bolt Flies | Titanium Nugget
$ws.Range("H74").Value = 10000
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

# Row 77: Heavy Metal Banned (L) | Titanium Nugget
$ws.Range("H77").Value = 10000
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value = 2010.4166
$ws.Range("I132").Value = 1612.5
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 4837.5
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -2307.5
$ws.Range("N132").Value = -17060

# Row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Range("H136").Value = 3522
$ws.Range("I136").Value = 2044
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 6132
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -3582
$ws.Range("N136").Value = -20100


$ws = $wb.Worksheets.Item("BSM")
# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 12118.412
$ws.Range("I134").Value = 14438.667
$ws.Range("K134").Value = 43316.001
$ws.Range("M134").Value = -40781.001

# Row 137: Dagger Swagger | Cobalt Tungsten Khukuri
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# Row 139: Maul Me | Titanium Gold Maul
$ws.Range("H139").Value = 40000
$ws.Range("J139").Value = 40000
$ws.Range("L139").Value = 40000
$ws.Range("N139").Value = -50280


$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 3270.375
$ws.Range("I31").Value = 1276.6428
$ws.Range("K31").Value = 1276.6428
$ws.Range("M31").Value = -981.6428000000001

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 3270.375
$ws.Range("I34").Value = 1276.6428
$ws.Range("K34").Value = 1276.6428
$ws.Range("M34").Value = -1074.6428

# Row 62: Splinter in the Sewers | Cedar Lumber
$ws.Range("H62").Value = 4980
$ws.Range("I62").Value = 4833.3335
$ws.Range("J62").Value = 5200
$ws.Range("K62").Value = 4833.3335
$ws.Range("L62").Value = 5200
$ws.Range("M62").Value = -4209.3335
$ws.Range("N62").Value = -6448

# Row 65: The Lumber of Their Discontent (L) | Cedar Lumber
$ws.Range("H65").Value = 4980
$ws.Range("I65").Value = 4833.3335
$ws.Range("J65").Value = 5200
$ws.Range("K65").Value = 24166.6675
$ws.Range("L65").Value = 26000
$ws.Range("M65").Value = -21046.6675
$ws.Range("N65").Value = -32240

# Row 134: Wood You Be Quiet | Ceiba Lumber
$ws.Range("H134").Value = 990.5333
$ws.Range("I134").Value = 949.61536
$ws.Range("K134").Value = 2848.84608
$ws.Range("M134").Value = -313.8460800000003


$ws = $wb.Worksheets.Item("CUL")
# Row 22: A Total Nut Job | Walnut Bread
$ws.Range("H22").Value = 3999.6667
$ws.Range("J22").Value = 4499.5
$ws.Range("L22").Value = 13498.5
$ws.Range("N22").Value = -13836.5

# Row 27: Brain Food | Walnut Bread
$ws.Range("H27").Value = 3999.6667
$ws.Range("J27").Value = 4499.5
$ws.Range("L27").Value = 13498.5
$ws.Range("N27").Value = -13702.5

# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Range("H131").Value = 12841661
$ws.Range("J131").Value = 26468.773
$ws.Range("L131").Value = 79406.319
$ws.Range("N131").Value = -89486.319


$ws = $wb.Worksheets.Item("GSM")
# Row 102: Put the Metal to the Peddle | Durium Ingot
$ws.Range("H102").Value = 1428.0667
$ws.Range("I102").Value = 1164.8889
$ws.Range("K102").Value = 1164.8889
$ws.Range("M102").Value = 457.1111000000001

# Row 113: Copious Crystal Cannons | Manasilver Nugget
$ws.Range("H113").Value = 1323.1428
$ws.Range("J113").Value = 1489.875
$ws.Range("L113").Value = 1489.875
$ws.Range("N113").Value = -5829.875

# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 3237
$ws.Range("I132").Value = 2917.7273
$ws.Range("J132").Value = 6749
$ws.Range("K132").Value = 8753.1819
$ws.Range("L132").Value = 20247
$ws.Range("M132").Value = -6223.1819
$ws.Range("N132").Value = -25307


$ws = $wb.Worksheets.Item("LTW")
# Row 36: Campaign in the Membrane | Toadskin Jacket
$ws.Range("H36").Value = 27715
$ws.Range("J36").Value = 27715
$ws.Range("L36").Value = 27715
$ws.Range("N36").Value = -28839

# Row 40: Best Served Toad | Toad Leather
$ws.Range("H40").Value = 8768.412
$ws.Range("J40").Value = 9782.667
$ws.Range("L40").Value = 9782.667
$ws.Range("N40").Value = -10054.667

# Row 46: Supply Side Logic | Boar Leather
$ws.Range("H46").Value = 1837.5
$ws.Range("J46").Value = 1837.5
$ws.Range("L46").Value = 1837.5
$ws.Range("N46").Value = -2213.5

# Row 82: Trainin' the Neck | Dragon Leather
$ws.Range("H82").Value = 4185
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()

# Row 85: Training Is Only Skintight (L) | Dragon Leather
$ws.Range("H85").Value = 4185
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()

# Row 122: Hell on Leather | Gaja Leather
$ws.Range("H122").Value = 6258.28
$ws.Range("I122").Value = 4550.1816
$ws.Range("J122").Value = 7600.357
$ws.Range("K122").Value = 13650.5448
$ws.Range("L122").Value = 22801.071
$ws.Range("M122").Value = -11200.5448
$ws.Range("N122").Value = -27701.071


$ws = $wb.Worksheets.Item("WVR")
# Row 13: Time for Acton | Hempen Acton
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()

# Row 81: Where the Dragonflies, the Net Catches | Crawler Silk
$ws.Range("H81").Value = 499.5
$ws.Range("I81").Value = 499
$ws.Range("K81").Value = 998
$ws.Range("M81").Value = 63

# Row 84: To Kill a Dragon on Nameday (L) | Crawler Silk
$ws.Range("H84").Value = 499.5
$ws.Range("I84").Value = 499
$ws.Range("K84").Value = 4990
$ws.Range("M84").Value = 314

# Row 95: Duress Rehearsal | Ruby Cotton Fingerless Gloves of Casting
$ws.Range("H95").Value = 48670
$ws.Range("J95").Value = 48670
$ws.Range("L95").Value = 48670
$ws.Range("N95").Value = -54162

# Row 122: Heavy Armoire | Dark Hempen Cloth
$ws.Range("H122").Value = 17188.344
$ws.Range("I122").Value = 23149.783
$ws.Range("J122").Value = 1953.5555
$ws.Range("K122").Value = 69449.349
$ws.Range("L122").Value = 5860.666499999999
$ws.Range("M122").Value = -66999.349
$ws.Range("N122").Value = -10760.6665

# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 4387.75
$ws.Range("I132").Value = 1663.125
$ws.Range("J132").Value = 7112.375
$ws.Range("K132").Value = 4989.375
$ws.Range("L132").Value = 21337.125
$ws.Range("M132").Value = -2459.375
$ws.Range("N132").Value = -26397.125

